$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23
$ws.Range("A23").Value = 112026957
$ws.Range("B23").Value = 85313
$ws.Range("C23").Value = 'Ovaliderad'
$ws.Range("D23").Value = 'NT'
$ws.Range("E23").Value = 3739
$ws.Range("F23").Value = 'Persiljespindling'
$ws.Range("G23").Value = 'Cortinarius sulfurinus'
$ws.Range("H23").Value = 'Quél.'
$ws.Range("P23").Value = 'Bye kalkbarrskogs naturreservat (Bye kalkbarrskogs naturreservat), Jmt'
$ws.Range("Q23").Value = 485421.0690379958
$ws.Range("R23").Value = 6996666.302800749
$ws.Range("S23").Value = 25
$ws.Range("T23").Value = 'Jämtland'
$ws.Range("U23").Value = 'Östersund'
$ws.Range("V23").Value = 'Jämtland'
$ws.Range("W23").Value = 'Marieby'
$ws.Range("Y23").Value = "'2023-09-11"
$ws.Range("Z23").Value = '00:00'
$ws.Range("AA23").Value = "'2023-09-11"
$ws.Range("AB23").Value = '00:00'
$ws.Range("AD23").Value = $false
$ws.Range("AE23").Value = $false
$ws.Range("AG23").Value = $false
$ws.Range("AW23").Value = 'Rashid Kadhim'
$ws.Range("AX23").Value = 'Rashid Kadhim'

# Row 24
$ws.Range("A24").Value = 112026905
$ws.Range("B24").Value = 90651
$ws.Range("C24").Value = 'Ovaliderad'
$ws.Range("D24").Value = 'NT'
$ws.Range("E24").Value = 1968
$ws.Range("F24").Value = 'Grantaggsvamp'
$ws.Range("G24").Value = 'Bankera violascens'
$ws.Range("H24").Value = '(Alb. & Schwein. : Fr.) Pouzar'
$ws.Range("P24").Value = 'Bye kalkbarrskogs naturreservat (Bye kalkbarrskogs naturreservat), Jmt'
$ws.Range("Q24").Value = 485427.0354010577
$ws.Range("R24").Value = 6996682.125945764
$ws.Range("S24").Value = 25
$ws.Range("T24").Value = 'Jämtland'
$ws.Range("U24").Value = 'Östersund'
$ws.Range("V24").Value = 'Jämtland'
$ws.Range("W24").Value = 'Marieby'
$ws.Range("Y24").Value = "'2023-09-11"
$ws.Range("Z24").Value = '00:00'
$ws.Range("AA24").Value = "'2023-09-11"
$ws.Range("AB24").Value = '00:00'
$ws.Range("AD24").Value = $false
$ws.Range("AE24").Value = $false
$ws.Range("AG24").Value = $false
$ws.Range("AW24").Value = 'Rashid Kadhim'
$ws.Range("AX24").Value = 'Rashid Kadhim'
